# Update "Förändrad" (Changed) date column (C) from 45326 (2024-02-04)
# to 45327 (2024-02-05) for all data rows (2-28) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45327
}
